$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a far-away scratch row as temp holding space for rotations so we
# never clobber a row before it has been read.
$scratchRow = 300

# --- Rotation 1: rows 142, 143, 145 (columns B:AC only; column A keeps the
#     sequential row index and is left untouched) ---
# New 142 <- old 145, new 143 <- old 142, new 145 <- old 143
$ws.Range("B${scratchRow}:AC${scratchRow}").Value2 = $ws.Range("B142:AC142").Value2
$ws.Range("B142:AC142").Value2 = $ws.Range("B145:AC145").Value2
$ws.Range("B145:AC145").Value2 = $ws.Range("B143:AC143").Value2
$ws.Range("B143:AC143").Value2 = $ws.Range("B${scratchRow}:AC${scratchRow}").Value2
$ws.Range("B${scratchRow}:AC${scratchRow}").ClearContents()

# --- Rotation 2: rows 148, 149, 150 (columns B:AC only) ---
# New 148 <- old 149, new 149 <- old 150, new 150 <- old 148
$ws.Range("B${scratchRow}:AC${scratchRow}").Value2 = $ws.Range("B148:AC148").Value2
$ws.Range("B148:AC148").Value2 = $ws.Range("B149:AC149").Value2
$ws.Range("B149:AC149").Value2 = $ws.Range("B150:AC150").Value2
$ws.Range("B150:AC150").Value2 = $ws.Range("B${scratchRow}:AC${scratchRow}").Value2
$ws.Range("B${scratchRow}:AC${scratchRow}").ClearContents()

# --- Swap: rows 211, 212 (columns B:AC only) ---
$ws.Range("B${scratchRow}:AC${scratchRow}").Value2 = $ws.Range("B211:AC211").Value2
$ws.Range("B211:AC211").Value2 = $ws.Range("B212:AC212").Value2
$ws.Range("B212:AC212").Value2 = $ws.Range("B${scratchRow}:AC${scratchRow}").Value2
$ws.Range("B${scratchRow}:AC${scratchRow}").ClearContents()

# --- Append new row 216, copying formatting (styles) from row 215 first so
#     the bold/border style on col A and the date format on col E carry
#     over exactly like every other data row. ---
$ws.Range("A215:AC215").Copy()
$ws.Range("A216:AC216").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A216").Value2 = 214
$ws.Range("B216").Value2 = 8039398
$ws.Range("C216").Value2 = "Bolivia Primera División"
$ws.Range("D216").Value2 = "Bolivia Apertura"
$ws.Range("E216").Value2 = 45389.85416666666
$ws.Range("F216").Value2 = "Real Tomayapo"
$ws.Range("G216").Value2 = "Real Santa Cruz"
$ws.Range("H216").Value2 = 2
$ws.Range("I216").Value2 = 1
$ws.Range("J216").Value2 = "H"
$ws.Range("K216").Value2 = 1.533
$ws.Range("L216").Value2 = 3.75
$ws.Range("M216").Value2 = 5.5
$ws.Range("N216").Value2 = 1.444
$ws.Range("O216").Value2 = 4.333
$ws.Range("P216").Value2 = 7.5
$ws.Range("Q216").Value2 = -1.25
$ws.Range("R216").Value2 = 2
$ws.Range("S216").Value2 = 1.8
$ws.Range("T216").Value2 = 2.75
$ws.Range("U216").Value2 = 1.975
$ws.Range("V216").Value2 = 1.825
$ws.Range("W216").Value2 = 0.444
$ws.Range("X216").Value2 = -1
$ws.Range("Y216").Value2 = -1
$ws.Range("Z216").Value2 = -0.5
$ws.Range("AA216").Value2 = 0.4
$ws.Range("AB216").Value2 = 0.4875
$ws.Range("AC216").Value2 = -0.5

Write-Host "edit complete"
